$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Val)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

Set-TextValue "D2" "29.240.12"
Set-TextValue "E2" "  +1.20%  "

Set-TextValue "D3" "1.838.50"
Set-TextValue "E3" "  +0.45%  "

Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  +0.06%  "

Set-TextValue "D5" "242.62"
Set-TextValue "E5" "  -0.91%  "

Set-TextValue "D6" "0.6864"
Set-TextValue "E6" "  -0.87%  "

Set-TextValue "D7" "1.001"
Set-TextValue "E7" "  +0.13%  "

Set-TextValue "D8" "0.3025"
Set-TextValue "E8" "  -0.40%  "

Set-TextValue "D9" "0.07509"
Set-TextValue "E9" "  -2.03%  "

Set-TextValue "D10" "23.25"
Set-TextValue "E10" "  -0.02%  "

Set-TextValue "D11" "0.07678"
Set-TextValue "E11" "  -1.66%  "

Set-TextValue "D12" "1.843.54"
Set-TextValue "E12" "  +0.65%  "

Set-TextValue "D13" "5.082"
Set-TextValue "E13" "  -0.16%  "

Set-TextValue "D14" "0.6860"
Set-TextValue "E14" "  +0.77%  "

Set-TextValue "D15" "88.46"
Set-TextValue "E15" "  -4.84%  "

Set-TextValue "D16" "6.257"
Set-TextValue "E16" "  -4.08%  "

Set-TextValue "D17" "29.293.65"
Set-TextValue "E17" "  +1.30%  "

Set-TextValue "D18" "0.000008203"
Set-TextValue "E18" "  -0.35%  "

Set-TextValue "D19" "2.087.83"
Set-TextValue "E19" "  +0.53%  "

Set-TextValue "D20" "231.94"
Set-TextValue "E20" "  -3.14%  "

Set-TextValue "D21" "12.59"
Set-TextValue "E21" "  -0.48%  "

Set-TextValue "D22" "0.9990"
Set-TextValue "E22" "  -0.12%  "

Set-TextValue "D23" "7.446"
Set-TextValue "E23" "  +0.01%  "

Set-TextValue "D24" "0.9994"
Set-TextValue "E24" "  -0.09%  "

Set-TextValue "D25" "0.1460"
Set-TextValue "E25" "  -2.51%  "

Set-TextValue "D26" "159.98"
Set-TextValue "E26" "  +1.30%  "

Set-TextValue "D27" "8.827"
Set-TextValue "E27" "  +1.15%  "

Set-TextValue "D28" "18.09"

Set-TextValue "D29" "1.517"
Set-TextValue "E29" "  -1.45%  "

Set-TextValue "D30" "4.266"
Set-TextValue "E30" "  +1.01%  "

Set-TextValue "D31" "4.144"
Set-TextValue "E31" "  +0.31%  "

Set-TextValue "D32" "1.208"
Set-TextValue "E32" "  +1.39%  "

Set-TextValue "D33" "0.05139"
Set-TextValue "E33" "  +0.60%  "

Set-TextValue "D34" "0.7686"
Set-TextValue "E34" "  -0.74%  "

Set-TextValue "D35" "1.837"
Set-TextValue "E35" "  -0.52%  "

Set-TextValue "D36" "1.137"
Set-TextValue "E36" "  -0.21%  "

Set-TextValue "D37" "2.672"
Set-TextValue "E37" "  -0.78%  "

Set-TextValue "D38" "1.305.69"
Set-TextValue "E38" "  +2.50%  "

Set-TextValue "E39" "  -0.52%  "

Set-TextValue "D40" "2.704"
Set-TextValue "E40" "  +0.26%  "

Set-TextValue "D41" "0.9430"
Set-TextValue "E41" "  -1.11%  "

Set-TextValue "B42" "Quant"
Set-TextValue "C42" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D42" "104.86"
Set-TextValue "E42" "  -1.77%  "

Set-TextValue "B43" "FraxShare"
Set-TextValue "C43" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "5.792"
Set-TextValue "E43" "  -5.75%  "

Set-TextValue "D44" "0.9988"
Set-TextValue "E44" "  -0.09%  "

Set-TextValue "D45" "65.42"
Set-TextValue "E45" "  +2.66%  "

Set-TextValue "D46" "9.628"
Set-TextValue "E46" "  -0.42%  "

Set-TextValue "D47" "1.987.69"
Set-TextValue "E47" "  +0.65%  "

Set-TextValue "B48" "BabyDogeCoin"
Set-TextValue "C48" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D48" "0.00000000123"
Set-TextValue "E48" "  +0.25%  "

Set-TextValue "B49" "Mantle"
Set-TextValue "C49" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D49" "0.5178"
Set-TextValue "E49" "  +0.26%  "

Set-TextValue "B50" "RenderToken"
Set-TextValue "C50" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D50" "1.771"
Set-TextValue "E50" "  +1.32%  "

Set-TextValue "D51" "0.05916"
Set-TextValue "E51" "  +0.83%  "
